# Update for first draft
# - Adds mean/std column pairs for each alt-horizon (previously a single column each)
# - Renames algorithm 'CART' -> 'DTREE'
# - Drops the 'NB' algorithm row entirely
# - Refreshes all numeric results with the new mean/std figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop the NB row (old row 9) first so everything below targets the
#    final 8-row (1 header + 7 data rows) layout.
$ws.Rows(9).Delete()

# 2. Make room for the 5 new std columns (H:L): copy the style of the
#    last existing header cell (G1) onto the new header cells.
$ws.Range("G1").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)

# 3. Header row text
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "One Year Alt mean"
$ws.Range("D1").Value = "One Year Alt std"
$ws.Range("E1").Value = "Two Year Alt mean"
$ws.Range("F1").Value = "Two Year Alt std"
$ws.Range("G1").Value = "Three Year Alt mean"
$ws.Range("H1").Value = "Three Year Alt std"
$ws.Range("I1").Value = "Five Year Alt mean"
$ws.Range("J1").Value = "Five Year Alt std"
$ws.Range("K1").Value = "Ten Year Alt mean"
$ws.Range("L1").Value = "Ten Year Alt std"

# 4. Data rows
# row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.8731713447842481
$ws.Range("D2").Value = 0.01899606356873269
$ws.Range("E2").Value = 0.8617760617760618
$ws.Range("F2").Value = 0.02925680658011167
$ws.Range("G2").Value = 0.850224752771951
$ws.Range("H2").Value = 0.02799383701769759
$ws.Range("I2").Value = 0.8491642441860465
$ws.Range("J2").Value = 0.04516712680155505
$ws.Range("K2").Value = 0.8429396000824573
$ws.Range("L2").Value = 0.04108206277159485

# row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.8789694176790951
$ws.Range("D3").Value = 0.01754735420855202
$ws.Range("E3").Value = 0.8665379665379666
$ws.Range("F3").Value = 0.02991401788066698
$ws.Range("G3").Value = 0.8643492158625513
$ws.Range("H3").Value = 0.02712078407341961
$ws.Range("I3").Value = 0.8647165697674419
$ws.Range("J3").Value = 0.04245675292791819
$ws.Range("K3").Value = 0.8733560090702948
$ws.Range("L3").Value = 0.03897908729589632

# row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.9443359865940512
$ws.Range("D4").Value = 0.01426937811627454
$ws.Range("E4").Value = 0.9376539805111234
$ws.Range("F4").Value = 0.02514167600603162
$ws.Range("G4").Value = 0.9420837079212866
$ws.Range("H4").Value = 0.01771964872633858
$ws.Range("I4").Value = 0.9479530038759689
$ws.Range("J4").Value = 0.02133556399830654
$ws.Range("K4").Value = 0.934126984126984
$ws.Range("L4").Value = 0.03776084186998975

# row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.7825764558022622
$ws.Range("D5").Value = 0.03975526215377703
$ws.Range("E5").Value = 0.7723616473616474
$ws.Range("F5").Value = 0.03645867638459473
$ws.Range("G5").Value = 0.7660823094595944
$ws.Range("H5").Value = 0.02985175149280916
$ws.Range("I5").Value = 0.7846778100775194
$ws.Range("J5").Value = 0.05411492589390134
$ws.Range("K5").Value = 0.7709235209235209
$ws.Range("L5").Value = 0.04109769075720256

# row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.8737997486384582
$ws.Range("D6").Value = 0.02218946677135712
$ws.Range("E6").Value = 0.8644649751792608
$ws.Range("F6").Value = 0.03014209697164604
$ws.Range("G6").Value = 0.874952552192588
$ws.Range("H6").Value = 0.02702395223789894
$ws.Range("I6").Value = 0.8786942829457365
$ws.Range("J6").Value = 0.03074453846014354
$ws.Range("K6").Value = 0.8621830550401979
$ws.Range("L6").Value = 0.02696539004305422

# row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.8841600335148723
$ws.Range("D7").Value = 0.02482985738127715
$ws.Range("E7").Value = 0.8780336458907888
$ws.Range("F7").Value = 0.02761690091688106
$ws.Range("G7").Value = 0.8742333433223454
$ws.Range("H7").Value = 0.02467888942738332
$ws.Range("I7").Value = 0.8631540697674419
$ws.Range("J7").Value = 0.03378929774449718
$ws.Range("K7").Value = 0.8662646876932591
$ws.Range("L7").Value = 0.02962388607795836

# row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.8977503142019272
$ws.Range("D8").Value = 0.019032466470846
$ws.Range("E8").Value = 0.8990255561684133
$ws.Range("F8").Value = 0.02552469984636539
$ws.Range("G8").Value = 0.8954599940065927
$ws.Range("H8").Value = 0.02203756334358568
$ws.Range("I8").Value = 0.906704215116279
$ws.Range("J8").Value = 0.02665763885094654
$ws.Range("K8").Value = 0.8905792620078333
$ws.Range("L8").Value = 0.02551332518819665

